$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 276.27274
$ws.Range("I9").Value = 188.375
$ws.Range("K9").Value = 188.375
$ws.Range("M9").Value = -19.375
$ws.Range("H18").Value = 328.2353
$ws.Range("I18").Value = 328.2353
$ws.Range("K18").Value = 328.2353
$ws.Range("M18").Value = -44.2353
$ws.Range("H19").Value = 1629.8235
$ws.Range("I19").Value = 1965.9231
$ws.Range("K19").Value = 1965.9231
$ws.Range("M19").Value = -1790.9231
$ws.Range("H33").Value = 731949.9
$ws.Range("I33").Value = 1287598.1
$ws.Range("J33").Value = 2661.625
$ws.Range("K33").Value = 1287598.1
$ws.Range("L33").Value = 2661.625
$ws.Range("M33").Value = -1287369.1
$ws.Range("N33").Value = -3119.625
$ws.Range("H40").Value = 2000.0385
$ws.Range("I40").Value = 1999.75
$ws.Range("J40").Value = 2000.091
$ws.Range("K40").Value = 1999.75
$ws.Range("L40").Value = 2000.091
$ws.Range("M40").Value = -1824.75
$ws.Range("N40").Value = -2350.091
$ws.Range("H43").Value = 2163.2273
$ws.Range("I43").Value = 1499.8422
$ws.Range("J43").Value = 6364.6665
$ws.Range("K43").Value = 1499.8422
$ws.Range("L43").Value = 6364.6665
$ws.Range("M43").Value = -1430.8422
$ws.Range("N43").Value = -6502.6665
$ws.Range("H47").Value = 17500
$ws.Range("J47").Value = 17500
$ws.Range("L47").Value = 17500
$ws.Range("N47").Value = -19444
$ws.Range("H55").Value = 193.76923
$ws.Range("I55").Value = 157.88889
$ws.Range("K55").Value = 157.88889
$ws.Range("M55").Value = 56.11111
$ws.Range("H62").Value = 9279.637000000001
$ws.Range("I62").Value = 9563
$ws.Range("J62").Value = 8004.5
$ws.Range("K62").Value = 9563
$ws.Range("L62").Value = 8004.5
$ws.Range("M62").Value = -8939
$ws.Range("N62").Value = -9252.5
$ws.Range("H65").Value = 9279.637000000001
$ws.Range("I65").Value = 9563
$ws.Range("J65").Value = 8004.5
$ws.Range("K65").Value = 47815
$ws.Range("L65").Value = 40022.5
$ws.Range("M65").Value = -44695
$ws.Range("N65").Value = -46262.5
$ws.Range("H70").Value = 3131.4546
$ws.Range("J70").Value = 3344.6
$ws.Range("L70").Value = 10033.8
$ws.Range("N70").Value = -10573.8
$ws.Range("H73").Value = 3131.4546
$ws.Range("J73").Value = 3344.6
$ws.Range("L73").Value = 10033.8
$ws.Range("N73").Value = -11905.8
$ws.Range("H80").Value = 477550.28
$ws.Range("J80").Value = 770979.9399999999
$ws.Range("L80").Value = 2312939.82
$ws.Range("N80").Value = -2314935.82
$ws.Range("H83").Value = 477550.28
$ws.Range("J83").Value = 770979.9399999999
$ws.Range("L83").Value = 6938819.459999999
$ws.Range("N83").Value = -6948803.459999999
$ws.Range("H86").Value = 1679.9375
$ws.Range("I86").Value = 875
$ws.Range("K86").Value = 875
$ws.Range("M86").Value = 248
$ws.Range("H89").Value = 1679.9375
$ws.Range("I89").Value = 875
$ws.Range("K89").Value = 4375
$ws.Range("M89").Value = 1241
$ws.Range("H92").Value = 557.8570999999999
$ws.Range("I92").Value = 536.1539
$ws.Range("J92").Value = 840
$ws.Range("K92").Value = 536.1539
$ws.Range("L92").Value = 840
$ws.Range("M92").Value = 711.8461
$ws.Range("N92").Value = -3336
$ws.Range("H100").Value = 38860.215
$ws.Range("I100").Value = 52411.3
$ws.Range("K100").Value = 52411.3
$ws.Range("M100").Value = -51870.3
$ws.Range("H101").Value = 1390
$ws.Range("I101").Value = 983.6667
$ws.Range("J101").Value = 1999.5
$ws.Range("K101").Value = 2951.0001
$ws.Range("L101").Value = 5998.5
$ws.Range("M101").Value = -1329.0001
$ws.Range("N101").Value = -9242.5
$ws.Range("H106").Value = 22128.348
$ws.Range("I106").Value = 9543.454
$ws.Range("K106").Value = 9543.454
$ws.Range("M106").Value = -8912.454
$ws.Range("H107").Value = 1417.8948
$ws.Range("I107").Value = 1021.625
$ws.Range("J107").Value = 3531.3333
$ws.Range("K107").Value = 1021.625
$ws.Range("L107").Value = 3531.3333
$ws.Range("M107").Value = 898.375
$ws.Range("N107").Value = -7371.3333
$ws.Range("H132").Value = 1322.6923
$ws.Range("I132").Value = 1259.5714
$ws.Range("K132").Value = 3778.7142
$ws.Range("M132").Value = -1248.7142
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 9112.554
$ws.Range("I137").Value = 4041.1086
$ws.Range("J137").Value = 17444.215
$ws.Range("K137").Value = 12123.3258
$ws.Range("L137").Value = 52332.645
$ws.Range("M137").Value = -9573.325800000001
$ws.Range("N137").Value = -57432.645
$ws.Range("H138").Value = 5214.3086
$ws.Range("I138").Value = 4227.394
$ws.Range("J138").Value = 6144.8286
$ws.Range("K138").Value = 12682.182
$ws.Range("L138").Value = 18434.4858
$ws.Range("M138").Value = -7542.182000000001
$ws.Range("N138").Value = -28714.4858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4761.492
$ws.Range("I32").Value = 2779
$ws.Range("K32").Value = 2779
$ws.Range("M32").Value = -2492
$ws.Range("H45").Value = 13205.6
$ws.Range("I45").Value = 14524.777
$ws.Range("K45").Value = 14524.777
$ws.Range("M45").Value = -14147.777
$ws.Range("H74").Value = 6222.8335
$ws.Range("I74").Value = 5257.6
$ws.Range("J74").Value = 6912.2856
$ws.Range("K74").Value = 5257.6
$ws.Range("L74").Value = 6912.2856
$ws.Range("M74").Value = -4383.6
$ws.Range("N74").Value = -8660.285599999999
$ws.Range("H77").Value = 6222.8335
$ws.Range("I77").Value = 5257.6
$ws.Range("J77").Value = 6912.2856
$ws.Range("K77").Value = 26288
$ws.Range("L77").Value = 34561.428
$ws.Range("M77").Value = -21920
$ws.Range("N77").Value = -43297.428
$ws.Range("H102").Value = 1524.4762
$ws.Range("I102").Value = 1503.2
$ws.Range("K102").Value = 1503.2
$ws.Range("M102").Value = 118.8
$ws.Range("H110").Value = 749
$ws.Range("I110").Value = 749
$ws.Range("J110").Value = 749
$ws.Range("K110").Value = 749
$ws.Range("L110").Value = 749
$ws.Range("M110").Value = 1296
$ws.Range("N110").Value = -4839
$ws.Range("H132").Value = 3951.0144
$ws.Range("I132").Value = 2694.6316
$ws.Range("K132").Value = 8083.8948
$ws.Range("M132").Value = -5553.8948

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 402.42307
$ws.Range("J80").Value = 427.9375
$ws.Range("L80").Value = 427.9375
$ws.Range("N80").Value = -2423.9375
$ws.Range("H83").Value = 402.42307
$ws.Range("J83").Value = 427.9375
$ws.Range("L83").Value = 2139.6875
$ws.Range("N83").Value = -12123.6875
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H99").Value = 3713.2666
$ws.Range("I99").Value = 2790.818
$ws.Range("K99").Value = 2790.818
$ws.Range("M99").Value = -1292.818
$ws.Range("H105").Value = 4575.8823
$ws.Range("I105").Value = 5024.25
$ws.Range("J105").Value = 3499.8
$ws.Range("K105").Value = 5024.25
$ws.Range("L105").Value = 3499.8
$ws.Range("M105").Value = -3277.25
$ws.Range("N105").Value = -6993.8
$ws.Range("H107").Value = 1249.3125
$ws.Range("I107").Value = 988.8889
$ws.Range("K107").Value = 988.8889
$ws.Range("M107").Value = 931.1111
$ws.Range("H134").Value = 12009.044
$ws.Range("I134").Value = 9883
$ws.Range("K134").Value = 29649
$ws.Range("M134").Value = -27114

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 328.1613
$ws.Range("I7").Value = 338.75
$ws.Range("K7").Value = 338.75
$ws.Range("M7").Value = -225.75
$ws.Range("H22").Value = 269.57895
$ws.Range("I22").Value = 239.125
$ws.Range("K22").Value = 239.125
$ws.Range("M22").Value = 110.875
$ws.Range("H31").Value = 70639.86
$ws.Range("I31").Value = 63600.562
$ws.Range("J31").Value = 79303.62
$ws.Range("K31").Value = 63600.562
$ws.Range("L31").Value = 79303.62
$ws.Range("M31").Value = -63305.562
$ws.Range("N31").Value = -79893.62
$ws.Range("H34").Value = 70639.86
$ws.Range("I34").Value = 63600.562
$ws.Range("J34").Value = 79303.62
$ws.Range("K34").Value = 63600.562
$ws.Range("L34").Value = 79303.62
$ws.Range("M34").Value = -63398.562
$ws.Range("N34").Value = -79707.62
$ws.Range("H42").Value = 25999
$ws.Range("J42").Value = 25999
$ws.Range("L42").Value = 25999
$ws.Range("N42").Value = -27185
$ws.Range("H50").Value = 49995
$ws.Range("J50").Value = 49995
$ws.Range("L50").Value = 49995
$ws.Range("N50").Value = -51245
$ws.Range("H58").Value = 21988.275
$ws.Range("I58").Value = 23552.355
$ws.Range("K58").Value = 23552.355
$ws.Range("M58").Value = -23349.355
$ws.Range("H59").Value = 34921.07
$ws.Range("I59").Value = 52499
$ws.Range("J59").Value = 31991.416
$ws.Range("K59").Value = 52499
$ws.Range("L59").Value = 31991.416
$ws.Range("M59").Value = -51354
$ws.Range("N59").Value = -34281.416
$ws.Range("H60").Value = 16373.125
$ws.Range("I60").Value = 2800
$ws.Range("J60").Value = 38995
$ws.Range("K60").Value = 2800
$ws.Range("L60").Value = 38995
$ws.Range("M60").Value = -2289
$ws.Range("N60").Value = -40017
$ws.Range("H94").Value = 3618
$ws.Range("J94").Value = 3618
$ws.Range("L94").Value = 3618
$ws.Range("N94").Value = -4520
$ws.Range("H105").Value = 2143.25
$ws.Range("I105").Value = 2118.2273
$ws.Range("K105").Value = 2118.2273
$ws.Range("M105").Value = -371.2273
$ws.Range("H107").Value = 1009.9
$ws.Range("I107").Value = 870.3333
$ws.Range("J107").Value = 1069.7142
$ws.Range("K107").Value = 870.3333
$ws.Range("L107").Value = 1069.7142
$ws.Range("M107").Value = 1049.6667
$ws.Range("N107").Value = -4909.7142
$ws.Range("H132").Value = 17951.453
$ws.Range("I132").Value = 11772.379
$ws.Range("J132").Value = 31735.54
$ws.Range("K132").Value = 35317.137
$ws.Range("L132").Value = 95206.62
$ws.Range("M132").Value = -32787.137
$ws.Range("N132").Value = -100266.62
$ws.Range("H136").Value = 21988.275
$ws.Range("I136").Value = 23552.355
$ws.Range("K136").Value = 70657.065
$ws.Range("M136").Value = -68107.065

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 142.44737
$ws.Range("I6").Value = 136.75757
$ws.Range("J6").Value = 180
$ws.Range("K6").Value = 410.27271
$ws.Range("L6").Value = 540
$ws.Range("M6").Value = -297.27271
$ws.Range("N6").Value = -766
$ws.Range("H8").Value = 701
$ws.Range("I8").Value = 701
$ws.Range("K8").Value = 2103
$ws.Range("M8").Value = -1964
$ws.Range("H11").Value = 374
$ws.Range("I11").Value = 359.29166
$ws.Range("K11").Value = 1077.87498
$ws.Range("M11").Value = -937.8749800000001
$ws.Range("H13").Value = 827.8333
$ws.Range("I13").Value = 83.5
$ws.Range("J13").Value = 1200
$ws.Range("K13").Value = 250.5
$ws.Range("L13").Value = 3600
$ws.Range("M13").Value = -82.5
$ws.Range("N13").Value = -3936
$ws.Range("H40").Value = 63.8
$ws.Range("I40").Value = 66.333336
$ws.Range("J40").Value = 62.714287
$ws.Range("K40").Value = 265.333344
$ws.Range("L40").Value = 250.857148
$ws.Range("M40").Value = -196.333344
$ws.Range("N40").Value = -388.857148
$ws.Range("H61").Value = 1273.5
$ws.Range("J61").Value = 1634.6666
$ws.Range("L61").Value = 4903.9998
$ws.Range("N61").Value = -5333.9998
$ws.Range("H94").Value = 6140
$ws.Range("I94").Value = 1700
$ws.Range("J94").Value = 7250
$ws.Range("K94").Value = 5100
$ws.Range("L94").Value = 21750
$ws.Range("M94").Value = -4424
$ws.Range("N94").Value = -23102
$ws.Range("H97").Value = 2284.5
$ws.Range("J97").Value = 5099.25
$ws.Range("L97").Value = 15297.75
$ws.Range("N97").Value = -16289.75
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 29999.166
$ws.Range("J15").Value = 29999.166
$ws.Range("L15").Value = 29999.166
$ws.Range("N15").Value = -30575.166
$ws.Range("H81").Value = 29999.166
$ws.Range("J81").Value = 29999.166
$ws.Range("L81").Value = 29999.166
$ws.Range("N81").Value = -31995.166
$ws.Range("H84").Value = 29999.166
$ws.Range("J84").Value = 29999.166
$ws.Range("L84").Value = 89997.49800000001
$ws.Range("N84").Value = -99981.49800000001
$ws.Range("H102").Value = 2276.5356
$ws.Range("I102").Value = 2154.9048
$ws.Range("K102").Value = 2154.9048
$ws.Range("M102").Value = -532.9047999999998
$ws.Range("H122").Value = 4565.143
$ws.Range("I122").Value = 4173.3335
$ws.Range("K122").Value = 12520.0005
$ws.Range("M122").Value = -10070.0005
$ws.Range("H126").Value = 2630.25
$ws.Range("I126").Value = 2740.7778
$ws.Range("K126").Value = 8222.3334
$ws.Range("M126").Value = -5752.3334
$ws.Range("H132").Value = 16043.111
$ws.Range("I132").Value = 21566.908
$ws.Range("J132").Value = 7362.857
$ws.Range("K132").Value = 64700.724
$ws.Range("L132").Value = 22088.571
$ws.Range("M132").Value = -62170.724
$ws.Range("N132").Value = -27148.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 1500
$ws.Range("I13").Value = 1500
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -1360
$ws.Range("N13").ClearContents()
$ws.Range("H16").Value = 1926.0667
$ws.Range("I16").Value = 1768.6154
$ws.Range("K16").Value = 1768.6154
$ws.Range("M16").Value = -1598.6154
$ws.Range("H40").Value = 5575.8335
$ws.Range("I40").Value = 5867
$ws.Range("J40").Value = 4993.5
$ws.Range("K40").Value = 5867
$ws.Range("L40").Value = 4993.5
$ws.Range("M40").Value = -5731
$ws.Range("N40").Value = -5265.5
$ws.Range("H46").Value = 1483.9354
$ws.Range("I46").Value = 1049.875
$ws.Range("J46").Value = 1634.9131
$ws.Range("K46").Value = 1049.875
$ws.Range("L46").Value = 1634.9131
$ws.Range("M46").Value = -861.875
$ws.Range("N46").Value = -2010.9131
$ws.Range("H55").Value = 132.21053
$ws.Range("I55").Value = 125.411766
$ws.Range("J55").Value = 190
$ws.Range("K55").Value = 125.411766
$ws.Range("L55").Value = 190
$ws.Range("M55").Value = 47.588234
$ws.Range("N55").Value = -536
$ws.Range("H61").Value = 2072
$ws.Range("I61").Value = 1828.4286
$ws.Range("K61").Value = 1828.4286
$ws.Range("M61").Value = -1626.4286
$ws.Range("H62").Value = 9800
$ws.Range("I62").Value = 9800
$ws.Range("K62").Value = 9800
$ws.Range("M62").Value = -9176
$ws.Range("H65").Value = 9800
$ws.Range("I65").Value = 9800
$ws.Range("K65").Value = 29400
$ws.Range("M65").Value = -26280
$ws.Range("H68").Value = 2829.138
$ws.Range("I68").Value = 2588.0952
$ws.Range("J68").Value = 3461.875
$ws.Range("K68").Value = 2588.0952
$ws.Range("L68").Value = 3461.875
$ws.Range("M68").Value = -1839.0952
$ws.Range("N68").Value = -4959.875
$ws.Range("H71").Value = 2829.138
$ws.Range("I71").Value = 2588.0952
$ws.Range("J71").Value = 3461.875
$ws.Range("K71").Value = 12940.476
$ws.Range("L71").Value = 17309.375
$ws.Range("M71").Value = -9196.476000000001
$ws.Range("N71").Value = -24797.375
$ws.Range("H113").Value = 2072
$ws.Range("I113").Value = 1828.4286
$ws.Range("K113").Value = 1828.4286
$ws.Range("M113").Value = 341.5714
$ws.Range("H132").Value = 5462.864
$ws.Range("I132").Value = 5171.1387
$ws.Range("K132").Value = 15513.4161
$ws.Range("M132").Value = -12983.4161
$ws.Range("H136").Value = 4948.2
$ws.Range("I136").Value = 5136.6
$ws.Range("J136").Value = 4759.8
$ws.Range("K136").Value = 15409.8
$ws.Range("L136").Value = 14279.4
$ws.Range("M136").Value = -12859.8
$ws.Range("N136").Value = -19379.4
$ws.Range("H139").Value = 26880
$ws.Range("I139").Value = 26880
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 26880
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -21740
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 493.75
$ws.Range("J6").Value = 493.75
$ws.Range("L6").Value = 493.75
$ws.Range("N6").Value = -723.75
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H74").Value = 23384.857
$ws.Range("J74").Value = 33332.5
$ws.Range("L74").Value = 33332.5
$ws.Range("N74").Value = -35204.5
$ws.Range("H77").Value = 23384.857
$ws.Range("J77").Value = 33332.5
$ws.Range("L77").Value = 99997.5
$ws.Range("N77").Value = -109357.5
$ws.Range("H100").Value = 203.625
$ws.Range("I100").Value = 117.583336
$ws.Range("K100").Value = 235.166672
$ws.Range("M100").Value = 305.833328
$ws.Range("H107").Value = 3099.1875
$ws.Range("I107").Value = 2371.5454
$ws.Range("K107").Value = 7114.6362
$ws.Range("M107").Value = -5194.6362
$ws.Range("H109").Value = 79371.375
$ws.Range("J109").Value = 79371.375
$ws.Range("L109").Value = 79371.375
$ws.Range("N109").Value = -82145.375
$ws.Range("H122").Value = 2899.1614
$ws.Range("I122").Value = 1718.8
$ws.Range("J122").Value = 5045.273
$ws.Range("K122").Value = 5156.4
$ws.Range("L122").Value = 15135.819
$ws.Range("M122").Value = -2706.4
$ws.Range("N122").Value = -20035.819
$ws.Range("H126").Value = 4931.7095
$ws.Range("I126").Value = 5301.1787
$ws.Range("K126").Value = 15903.5361
$ws.Range("M126").Value = -13433.5361
$ws.Range("H132").Value = 18046.61
$ws.Range("I132").Value = 9860.634
$ws.Range("J132").Value = 32639
$ws.Range("K132").Value = 29581.902
$ws.Range("L132").Value = 97917
$ws.Range("M132").Value = -27051.902
$ws.Range("N132").Value = -102977
$ws.Range("H136").Value = 1931.1428
$ws.Range("I136").Value = 1322.6957
$ws.Range("K136").Value = 3968.0871
$ws.Range("M136").Value = -1418.0871
$ws.Range("H139").Value = 58657.8
$ws.Range("I139").Value = 44999
$ws.Range("J139").Value = 62072.5
$ws.Range("K139").Value = 44999
$ws.Range("L139").Value = 62072.5
$ws.Range("M139").Value = -39859
$ws.Range("N139").Value = -72352.5
$ws.Range("H141").Value = 53772
$ws.Range("J141").Value = 53772
$ws.Range("L141").Value = 53772
$ws.Range("N141").Value = -64132
